$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.689.37'
$ws.Range('E2').Value = '  +2.85%  '
$ws.Range('D3').Value = '4.025.56'
$ws.Range('E3').Value = '  +2.16%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '525.15'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '149.18'
$ws.Range('E6').Value = '  +1.79%  '
$ws.Range('E7').Value = '  +0.70%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.742'
$ws.Range('E9').Value = '  +1.49%  '
$ws.Range('E10').Value = '  +1.59%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000343'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '46.09'
$ws.Range('E12').Value = '  +7.35%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.77'
$ws.Range('E13').Value = '  +2.99%  '
$ws.Range('D14').Value = '4.671.05'
$ws.Range('E14').Value = '  +2.25%  '
$ws.Range('D15').Value = '4.035.20'
$ws.Range('E15').Value = '  +2.27%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.45'
$ws.Range('E16').Value = '  +7.64%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.30'
$ws.Range('E17').Value = '  +1.50%  '
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('D20').Value = '71.620.20'
$ws.Range('E20').Value = '  +2.90%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '441.39'
$ws.Range('E21').Value = '  +1.60%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.60'
$ws.Range('E22').Value = '  +5.93%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '94.81'
$ws.Range('E23').Value = '  +6.96%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '14.38'
$ws.Range('E24').Value = '  -1.48%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.22'
$ws.Range('E25').Value = '  +2.58%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.06'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.24'
$ws.Range('E27').Value = '  +2.20%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '37.06'
$ws.Range('E28').Value = '  +0.79%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '13.63'
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '699.79'
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.131'
$ws.Range('E31').Value = '  +2.72%  '
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('E33').Value = '  +13.99%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '67.72'
$ws.Range('E34').Value = '  -1.42%  '
$ws.Range('D35').Value = '0.0₃0910'
$ws.Range('E35').Value = '  +4.29%  '
$ws.Range('E36').Value = '  -1.54%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '41.02'
$ws.Range('E37').Value = '  +0.93%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.158'
$ws.Range('E38').Value = '  +5.88%  '
$ws.Range('E39').Value = '  +18.94%  '
$ws.Range('E40').Value = '  +0.25%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0491'
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.13'
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.83'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.51'
$ws.Range('E45').Value = '  +2.89%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.147'
$ws.Range('E46').Value = '  +2.42%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.19'
$ws.Range('E47').Value = '  -2.47%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.20'
$ws.Range('E48').Value = '  +5.93%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.000277'
$ws.Range('E49').Value = '  +17.31%  '
$ws.Range('B50').Value = 'LidoDAOToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.38'
$ws.Range('E50').Value = '  +0.93%  '
$ws.Range('E51').Value = '  -0.84%  '
